$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph that follows the H1 title ---
$p2 = $d.Paragraphs.Item(2)
if ($p2.Range.Text -like "Meta description*") {
    $p2.Range.Delete() | Out-Null
}

# --- Step 2: replace the "Prompt: ..." text with the meta-description text,
#             keeping its italic formatting intact ---
$oldPrompt = 'Prompt: Create a feature image for "Football Mania Deluxe" that is in cartoon style and features a happy Maya warrior with glasses. The warrior should be shown in a football jersey and be holding a football in one hand while giving a thumbs-up sign with the other hand. The background should be a football field with fans cheering in the stands. Make the image lively and colorful, using bright colors to make it stand out.'
$newDescription = 'Find out about the special features, immersive environment, symbol design, and customizable gaming options of Football Mania Deluxe. Play free now.'
$null = $d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)

# --- Step 3: insert a new bold paragraph with the page title right before it ---
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($n)
$newPara = $pLast.Range.InsertParagraphBefore()

$n2 = $d.Paragraphs.Count
$pTitle = $d.Paragraphs.Item($n2 - 1)
$titleRange = $pTitle.Range

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Football Mania Deluxe Free Slot | Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $titleRange.InsertXML($xmlFrag)
